$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 91, shifting the existing rows 91..180 down to 92..181
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row with the new price-report record
$ws.Range("A91").Value = 3
$ws.Range("B91").Value = "Femacal de La Calera"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44790
$ws.Range("E91").Value = 5
$ws.Range("F91").Value = 100112026
$ws.Range("G91").Value = "Haba"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 90
$ws.Range("K91").Value = 16000
$ws.Range("L91").Value = 17000
$ws.Range("M91").Value = 16500
$ws.Range("N91").Value = "$/saco 25 kilos"
$ws.Range("O91").Value = "Provincia de Limarí"
$ws.Range("P91").Value = 660
$ws.Range("Q91").Value = 25
$ws.Range("R91").Value = "Hortaliza"
